# Trade #71 (momentum strategy, global Trade #99) closed at 2026-02-18 00:27:21
# and a new MarketMaking trade (#128) opened at 2026-02-18 00:27:15.
# Updates Summary, Strategy Status, All Trades and the per-strategy sheets
# (momentum / MarketMaking) to reflect the closed trade + the new open trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.26
$wsSummary.Range("B4").Value = 0.37
$wsSummary.Range("B5").Value = 0.07000000000000001
$wsSummary.Range("B6").Value = 99
$wsSummary.Range("B7").Value = 47
$wsSummary.Range("B9").Value = 47.47

# ---------------------------------------------------------------------------
# Sheet: Strategy Status (row 11 = momentum)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C11").Value = 99.26000000000001
$wsStatus.Range("D11").Value = 22
$wsStatus.Range("E11").Value = -0.74
$wsStatus.Range("F11").Value = -0.74
$wsStatus.Range("G11").Value = 22.73

# ---------------------------------------------------------------------------
# Sheet: All Trades - close out trade #99 (row 100, momentum)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Cells.Item(100, 7).Value = 0.71           # G100 Exit Price
$wsAll.Cells.Item(100, 8).Value = "CLOSED"        # H100 Status
$wsAll.Cells.Item(100, 9).Value = 9.2308          # I100 P&L %
$wsAll.Cells.Item(100, 10).Value = 0.06           # J100 P&L $
$wsAll.Cells.Item(100, 11).Value = 99.26000000000001  # K100 Capital After
$wsAll.Cells.Item(100, 12).Value = "early_exit"   # L100 Exit Reason
$wsAll.Cells.Item(100, 13).Value = 0.13           # M100 Duration (min)

# ---------------------------------------------------------------------------
# Sheet: momentum - same trade, local row 23
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Cells.Item(23, 7).Value = 0.71            # G23 Exit Price
$wsMomentum.Cells.Item(23, 8).Value = "CLOSED"         # H23 Status
$wsMomentum.Cells.Item(23, 9).Value = 9.2308           # I23 P&L %
$wsMomentum.Cells.Item(23, 10).Value = 0.06            # J23 P&L $
$wsMomentum.Cells.Item(23, 11).Value = 99.26000000000001  # K23 Capital After
$wsMomentum.Cells.Item(23, 16).Value = "early_exit"    # P23 Exit Reason
$wsMomentum.Cells.Item(23, 17).Value = 0.13            # Q23 Duration (min)

# ---------------------------------------------------------------------------
# Sheet: All Trades - append newly opened trade #128 (MarketMaking) as row 129
# ---------------------------------------------------------------------------
# Date/Time-looking text must be copied from an existing text cell rather than
# typed directly, otherwise Excel auto-converts "2026-02-18" into a date serial.
$wsAll.Cells.Item(129, 1).Value = 128
$wsAll.Range("B31").Copy($wsAll.Cells.Item(129, 2))   # B129 Date "2026-02-18"
$wsAll.Cells.Item(129, 3).Value = "00:27:15"          # C129 Time
$wsAll.Cells.Item(129, 4).Value = "MarketMaking"      # D129 Strategy
$wsAll.Cells.Item(129, 5).Value = "DOWN"              # E129 Side
$wsAll.Cells.Item(129, 6).Value = 0.65                # F129 Entry Price
$wsAll.Cells.Item(129, 8).Value = "OPEN"              # H129 Status
$wsAll.Cells.Item(129, 9).Value = 0                   # I129 P&L %
$wsAll.Cells.Item(129, 10).Value = 0                  # J129 P&L $
$wsAll.Cells.Item(129, 11).Value = 99.40967800952272  # K129 Capital After
$wsAll.Cells.Item(129, 13).Value = 0                  # M129 Duration (min)
$wsAll.Cells.Item(129, 14).Value = 0                  # N129 Entry Slippage
$wsAll.Cells.Item(129, 15).Value = 0                  # O129 Exit Slippage
$wsAll.Cells.Item(129, 16).Value = 0.65               # P129 Confidence
$wsAll.Cells.Item(129, 17).Value = "Wide spread capture: 392 bps vs avg 283 bps"  # Q129 Entry Reason

# ---------------------------------------------------------------------------
# Sheet: MarketMaking - append the same new trade as local row 49
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Cells.Item(49, 1).Value = 128
$wsMM.Range("B2").Copy($wsMM.Cells.Item(49, 2))       # B49 Date "2026-02-18"
$wsMM.Cells.Item(49, 3).Value = "00:27:15"            # C49 Time
$wsMM.Cells.Item(49, 4).Value = "MarketMaking"        # D49 Strategy
$wsMM.Cells.Item(49, 5).Value = "DOWN"                # E49 Side
$wsMM.Cells.Item(49, 6).Value = 0.65                  # F49 Entry Price
$wsMM.Cells.Item(49, 8).Value = "OPEN"                # H49 Status
$wsMM.Cells.Item(49, 9).Value = 0                     # I49 P&L %
$wsMM.Cells.Item(49, 10).Value = 0                    # J49 P&L $
$wsMM.Cells.Item(49, 11).Value = 99.40967800952272    # K49 Capital After
$wsMM.Cells.Item(49, 12).Value = 0                    # L49 Entry Slippage
$wsMM.Cells.Item(49, 13).Value = 0                    # M49 Exit Slippage
$wsMM.Cells.Item(49, 14).Value = 0.65                 # N49 Confidence
$wsMM.Cells.Item(49, 15).Value = "Wide spread capture: 392 bps vs avg 283 bps"  # O49 Entry Reason
$wsMM.Cells.Item(49, 17).Value = 0                    # Q49 Duration (min)
